$d = $word.ActiveDocument

# Locate the list item paragraph that ends with "Poolside Alarm by Code" -
# the new "Baby Fence" bullet goes directly after it and before
# "Pool School Instruction".
$targetIdx = 0
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*Poolside Alarm by Code*") {
        $targetIdx = $idx
        break
    }
}

if ($targetIdx -eq 0) {
    Write-Host "Could not locate insertion point paragraph."
} else {
    $target = $d.Paragraphs.Item($targetIdx)

    # Create a new paragraph right after it; it inherits the same list
    # (numId 54) / run formatting as its neighbours.
    $target.Range.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($targetIdx + 1)
    $newPara.Range.Text = "Install 50 L.F. of Baby Fence accessing pool area by code"

    Write-Host "Inserted:" $newPara.Range.Text
}
